$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Q1: header text change
$ws.Range("Q1").Value = "wtkappa.scale_trim"

# E2, F2: sign flips
$ws.Range("E2").Value = 0.02351246133036713
$ws.Range("F2").Value = -0.01081937260331701

# Q2: value update
$ws.Range("Q2").Value = 0.7808705382933534
